# Fixed variables and query errors in Bread from TC01 to TC30
#
# The "CasesTab" query stored in cell B2 of the "startup" sheet had a
# trailing clause referencing a `cohort` node/variable that is not part of
# the intended query (it duplicated logic that belongs to the FilesTab
# query and caused a query error). Remove the erroneous trailing
# ", coalesce(co.cohort_description, '') AS `Cohort`" clause so the query
# ends cleanly after the "Response to Treatment" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected query text, base64-encoded (UTF-8) to avoid any PowerShell
# escaping issues with backticks/newlines/quotes inside the Cypher text.
$b64 = "TUFUQ0ggKHM6c3R1ZHkpPC1bKl0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYykKV0hFUkUgZGVtby5icmVlZCBJTlsnRW5nbGlzaCBTZXR0ZXInXSAKTUFUQ0ggKGMpPC0tKGRpYWc6ZGlhZ25vc2lzKQpPUFRJT05BTCBNQVRDSCAoc2FtcDpzYW1wbGUpLS0+KGMpCk9QVElPTkFMIE1BVENIIChjbzpjb2hvcnQpPC1bKl0tKGMpCldJVEggRElTVElOQ1QgYywgcywgZGVtbywgZGlhZywgY28KUkVUVVJOICBjb2FsZXNjZShjLmNhc2VfaWQsICcnKSBBUyBgQ2FzZSBJRGAgLAogICAgICAgIGNvYWxlc2NlKHMuY2xpbmljYWxfc3R1ZHlfZGVzaWduYXRpb24sICcnKSBBUyBgU3R1ZHkgQ29kZWAgLAogICAgICAgIGNvYWxlc2NlKHMuY2xpbmljYWxfc3R1ZHlfdHlwZSwgJycpIEFTICBgU3R1ZHkgVHlwZWAsCiAgICAgICAgY29hbGVzY2UoZGVtby5icmVlZCwgJycpIEFTIEJyZWVkICwKICAgICAgICBjb2FsZXNjZShkaWFnLmRpc2Vhc2VfdGVybSwgJycpIEFTIERpYWdub3NpcyAsCiAgICAgICAgY29hbGVzY2UoZGlhZy5zdGFnZV9vZl9kaXNlYXNlLCAnJykgQVMgYFN0YWdlIG9mIERpc2Vhc2VgICwKICAgICAgICBjb2FsZXNjZShkZW1vLnBhdGllbnRfYWdlX2F0X2Vucm9sbG1lbnQsICcnKSBBUyBBZ2UgLAogICAgICAgIGNvYWxlc2NlKGRlbW8uc2V4LCAnJykgQVMgU2V4ICwKICAgICAgICBjb2FsZXNjZShkZW1vLm5ldXRlcmVkX2luZGljYXRvciwgJycpIEFTIGBOZXV0ZXJlZCBTdGF0dXNgLAogICAgICAgIGNvYWxlc2NlKGRlbW8ud2VpZ2h0LCAnJykgQVMgYFdlaWdodCAoa2cpYCwKICAgICAgICBjb2FsZXNjZShkaWFnLmJlc3RfcmVzcG9uc2UsICcnKSBBUyBgUmVzcG9uc2UgdG8gVHJlYXRtZW50YA=="
$bytes = [System.Convert]::FromBase64String($b64)
$newCasesQuery = [System.Text.Encoding]::UTF8.GetString($bytes)

$ws.Range("B2").Value = $newCasesQuery

# Reflect the author's final cursor/selection position (cell B2, scrolled
# to the top of the sheet) instead of the previous B4 / scrolled-down view.
$ws.Range("B2").Select() | Out-Null
